# Auto-generated script to apply Golem_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H12").Value = 1030.125
$ws.Range("I12").Value = 974.25
$ws.Range("J12").Value = 1086
$ws.Range("K12").Value = 974.25
$ws.Range("L12").Value = 1086
$ws.Range("M12").Value = -804.25
$ws.Range("N12").Value = -1426
$ws.Range("H58").Value = 48.333332
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H62").Value = 3920
$ws.Range("J62").Value = 4450
$ws.Range("L62").Value = 4450
$ws.Range("N62").Value = -5698
$ws.Range("H65").Value = 3920
$ws.Range("J65").Value = 4450
$ws.Range("L65").Value = 22250
$ws.Range("N65").Value = -28490
$ws.Range("H125").Value = 1003.5714
$ws.Range("I125").Value = 997.25
$ws.Range("K125").Value = 8975.25
$ws.Range("M125").Value = -6515.25
$ws.Range("H137").Value = 2922.5
$ws.Range("J137").Value = 4125
$ws.Range("L137").Value = 12375
$ws.Range("N137").Value = -17475

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1499.75
$ws.Range("I45").Value = 1599.6666
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 1599.6666
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -1222.6666
$ws.Range("N45").Value = -1954
$ws.Range("H55").Value = 50000
$ws.Range("I55").Value = 50000
$ws.Range("K55").Value = 50000
$ws.Range("M55").Value = -49685
$ws.Range("H61").Value = 3797.6667
$ws.Range("I61").Value = 3794
$ws.Range("K61").Value = 3794
$ws.Range("M61").Value = -3582
$ws.Range("H88").Value = 955
$ws.Range("J88").Value = 955
$ws.Range("L88").Value = 955
$ws.Range("N88").Value = -1767
$ws.Range("H91").Value = 955
$ws.Range("J91").Value = 955
$ws.Range("L91").Value = 955
$ws.Range("N91").Value = -3763
$ws.Range("H102").Value = 17501180
$ws.Range("I102").Value = 910379.2
$ws.Range("J102").Value = 200000000
$ws.Range("K102").Value = 910379.2
$ws.Range("L102").Value = 200000000
$ws.Range("M102").Value = -908757.2
$ws.Range("N102").Value = -200003244
$ws.Range("H136").Value = 3797.6667
$ws.Range("I136").Value = 3794
$ws.Range("K136").Value = 11382
$ws.Range("M136").Value = -8832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 52167.6
$ws.Range("I26").Value = 43612.668
$ws.Range("K26").Value = 43612.668
$ws.Range("M26").Value = -43320.668
$ws.Range("H47").Value = 199999
$ws.Range("J47").Value = 199999
$ws.Range("L47").Value = 199999
$ws.Range("N47").Value = -201039
$ws.Range("H94").Value = 101353.09
$ws.Range("I94").Value = 123652.664
$ws.Range("J94").Value = 1005
$ws.Range("K94").Value = 123652.664
$ws.Range("L94").Value = 1005
$ws.Range("M94").Value = -123201.664
$ws.Range("N94").Value = -1907
$ws.Range("H96").Value = 2357
$ws.Range("I96").Value = 2357
$ws.Range("K96").Value = 2357
$ws.Range("M96").Value = 389
$ws.Range("H99").Value = 3583.3333
$ws.Range("I99").Value = 6000
$ws.Range("K99").Value = 6000
$ws.Range("M99").Value = -4502
$ws.Range("H104").Value = 70000
$ws.Range("J104").Value = 70000
$ws.Range("L104").Value = 70000
$ws.Range("N104").Value = -76988
$ws.Range("H107").Value = 44779.445
$ws.Range("I107").Value = 44779.445
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 44779.445
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -42859.445
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 879
$ws.Range("I134").Value = 879
$ws.Range("K134").Value = 2637
$ws.Range("M134").Value = -102

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1389556.4
$ws.Range("I6").Value = 1851075.1
$ws.Range("K6").Value = 1851075.1
$ws.Range("M6").Value = -1850962.1
$ws.Range("H7").Value = 409
$ws.Range("I7").Value = 409
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 409
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -296
$ws.Range("N7").ClearContents()
$ws.Range("H31").Value = 6796
$ws.Range("I31").Value = 2903.6667
$ws.Range("J31").Value = 9131.4
$ws.Range("K31").Value = 2903.6667
$ws.Range("L31").Value = 9131.4
$ws.Range("M31").Value = -2608.6667
$ws.Range("N31").Value = -9721.4
$ws.Range("H34").Value = 6796
$ws.Range("I34").Value = 2903.6667
$ws.Range("J34").Value = 9131.4
$ws.Range("K34").Value = 2903.6667
$ws.Range("L34").Value = 9131.4
$ws.Range("M34").Value = -2701.6667
$ws.Range("N34").Value = -9535.4
$ws.Range("H62").Value = 4500
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3876
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4500
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -19380
$ws.Range("N65").ClearContents()
$ws.Range("H107").Value = 630.6
$ws.Range("I107").Value = 595.3333
$ws.Range("K107").Value = 595.3333
$ws.Range("M107").Value = 1324.6667
$ws.Range("H134").Value = 2854.5715
$ws.Range("I134").Value = 2995.5
$ws.Range("K134").Value = 8986.5
$ws.Range("M134").Value = -6451.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1964.5927
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 2001.6923
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 6005.0769
$ws.Range("M22").Value = -2831
$ws.Range("N22").Value = -6343.0769
$ws.Range("H27").Value = 1964.5927
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 2001.6923
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 6005.0769
$ws.Range("M27").Value = -2898
$ws.Range("N27").Value = -6209.0769
$ws.Range("H69").Value = 1495
$ws.Range("I69").Value = 1495
$ws.Range("K69").Value = 4485
$ws.Range("M69").Value = -3674
$ws.Range("H72").Value = 1495
$ws.Range("I72").Value = 1495
$ws.Range("K72").Value = 13455
$ws.Range("M72").Value = -9399
$ws.Range("H98").Value = 140
$ws.Range("I98").Value = 140
$ws.Range("K98").Value = 420
$ws.Range("M98").Value = 1078

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 17062.5
$ws.Range("J92").Value = 17062.5
$ws.Range("L92").Value = 17062.5
$ws.Range("N92").Value = -20806.5
$ws.Range("H98").Value = 8629
$ws.Range("J98").Value = 8629
$ws.Range("L98").Value = 8629
$ws.Range("N98").Value = -14619
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H113").Value = 884.8333
$ws.Range("I113").Value = 884.8333
$ws.Range("K113").Value = 884.8333
$ws.Range("M113").Value = 1285.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9125
$ws.Range("I7").Value = 8833.333000000001
$ws.Range("K7").Value = 8833.333000000001
$ws.Range("M7").Value = -8721.333000000001
$ws.Range("H122").Value = 5453.222
$ws.Range("I122").Value = 4288.5
$ws.Range("J122").Value = 5786
$ws.Range("K122").Value = 12865.5
$ws.Range("L122").Value = 17358
$ws.Range("M122").Value = -10415.5
$ws.Range("N122").Value = -22258
$ws.Range("H126").Value = 9125
$ws.Range("I126").Value = 8833.333000000001
$ws.Range("K126").Value = 26499.999
$ws.Range("M126").Value = -24029.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 26499.8
$ws.Range("I55").Value = 8499.5
$ws.Range("K55").Value = 8499.5
$ws.Range("M55").Value = -8222.5
